# Refresh the cryptocurrency Price (col D) and Volume(1h) (col E) figures for the
# Sun Jan 22 2023 data pull (GitHub Actions scheduled scrape).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column D/E hold numeric-looking values stored as literal text (e.g. "301.93",
# "-0.92%") so that trailing zeros / percent signs render exactly as scraped.
# A leading apostrophe tells Excel to keep the entry as text instead of coercing
# it to a Number, matching the source workbook's existing cell type.
$cellUpdates = [ordered]@{
    "D2" = "301.93"
    "E2" = "-0.92%"
    "D3" = "37.45"
    "E3" = "7.19%"
    "D4" = "5.010"
    "E4" = "-2.64%"
    "D5" = "0.07840"
    "E5" = "0.71%"
    "D6" = "2.221"
    "E6" = "-6.35%"
    "E7" = "-0.32%"
    "D8" = "4.020"
    "D9" = "0.9083"
    "E9" = "-1.82%"
    "D10" = "0.09604"
    "E10" = "-2.25%"
    "D11" = "0.1894"
    "E11" = "4.91%"
    "D12" = "0.08489"
    "E12" = "-2.80%"
    "D13" = "0.03522"
    "E13" = "6.36%"
    "D14" = "0.09954"
    "E14" = "0.35%"
    "D15" = "0.001479"
    "E15" = "-1.03%"
    "D16" = "0.005644"
    "E16" = "-0.50%"
    "E17" = "-0.14%"
    "E18" = "-4.18%"
    "E19" = "2.82%"
    "E20" = "0.00%"
    "D21" = "4.774"
    "E21" = "9.60%"
    "D23" = "0.04647"
    "E23" = "1.73%"
    "E24" = "1.09%"
    "E25" = "0.03%"
    "E26" = "-7.57%"
    "D27" = "0.0004750"
    "E27" = "28.49%"
    "D39" = "0.01759"
    "E39" = "-2.17%"
    "D40" = "0.04724"
    "E40" = "-0.73%"
    "D41" = "0.007856"
    "E41" = "1.76%"
    "D42" = "0.1392"
    "E42" = "-1.40%"
    "D43" = "0.007663"
    "E43" = "5.49%"
    "D44" = "0.002162"
    "E44" = "-2.47%"
    "D45" = "0.009881"
    "E45" = "3.65%"
    "D46" = "0.00006065"
    "E46" = "-0.81%"
    "E47" = "0.15%"
    "D49" = "0.002690"
    "E49" = "34.61%"
    "E50" = "0.15%"
    "E51" = "0.15%"
}

foreach ($addr in $cellUpdates.Keys) {
    $ws.Range($addr).Value = "'" + $cellUpdates[$addr]
}
